# Deploying to gh-pages from @ NIH-NCPI/ncpi-fhir-ig-2@b701e861ff4aea87f49ab6a6b6da8d47ed8dfde7
# Regenerated content-version StructureDefinition export: refreshed the
# build Date, retargeted FHIR Version back to 4.0.1 (R4) and fixed the
# FHIR-R4-related text that had drifted to R4B/4.3.0 wording.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value  = "2025-06-13T15:45:04+00:00"   # Date
$meta.Range("B15").Value = "4.0.1"                         # FHIR Version

# ---- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 = "Extension" (base element) - Constraint(s) column AJ:
# drop the "unless an empty Parameters resource ... or $this is Parameters" clause
$elements.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 3 = "Extension.id" - Type(s) column K: "id" -> "string"
$elements.Range("K3").Value = "string`n"

# Row 6 = "Extension.value[x]" - Definition column M: R4B -> R4 doc link
$elements.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
